$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.987.05"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").Value = "1.896.61"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.8332"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.07%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "241.74"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3279"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.99%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "26.52"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07040"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "

# Row 11
$ws.Range("E11").Value = "  +0.48%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.7598"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "

# Row 13
$ws.Range("D13").Value = "1.898.66"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.242"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "92.17"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "

# Row 16
$ws.Range("D16").Value = "29.990.97"
$ws.Range("E16").Value = "  -0.36%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.07"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "5.861"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.47%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "244.42"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007755"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").Value = "2.146.38"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.963"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1729"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +22.79%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.239"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.75%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "165.86"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

# Row 28
$ws.Range("E28").Value = "  -0.67%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.093"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.55%  "

# Row 30
$ws.Range("E30").Value = "  -2.18%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.517"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.05881"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +9.49%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.279"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.95%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.073"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "

# Row 35
$ws.Range("E35").Value = "  -0.33%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7293"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.68%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.726"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01915"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.774"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.4426"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "72.31"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.8566"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.93%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.849"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -5.27%  "

# Row 44
$ws.Range("E44").Value = "  +0.05%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.886"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "101.78"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "7.547"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.792"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "997.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.30%  "

# Row 50
$ws.Range("D50").Value = "2.046.12"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "35.82"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
